$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $NewValue)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $NewValue
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '93.429.78'
Set-TextValue $ws.Range('E2') '  -4.42%  '
Set-TextValue $ws.Range('D3') '3.435.10'
Set-TextValue $ws.Range('E3') '  +2.25%  '
Set-TextValue $ws.Range('E4') '  +0.00%  '
Set-TextValue $ws.Range('D5') '235.74'
Set-TextValue $ws.Range('E5') '  -6.88%  '
Set-TextValue $ws.Range('D6') '637.67'
Set-TextValue $ws.Range('E6') '  -3.51%  '
Set-TextValue $ws.Range('E7') '  -0.47%  '
Set-TextValue $ws.Range('D8') '0.394'
Set-TextValue $ws.Range('E8') '  -7.94%  '
Set-TextValue $ws.Range('E9') '  +0.18%  '
Set-TextValue $ws.Range('E10') '  -6.90%  '
Set-TextValue $ws.Range('D11') '3.436.24'
Set-TextValue $ws.Range('E11') '  +2.35%  '
Set-TextValue $ws.Range('D12') '41.68'
Set-TextValue $ws.Range('E12') '  -0.87%  '
Set-TextValue $ws.Range('E13') '  -5.84%  '
Set-TextValue $ws.Range('E14') '  +0.22%  '
Set-TextValue $ws.Range('B15') 'WrappedBTC'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D15') '93.420.38'
Set-TextValue $ws.Range('E15') '  -4.25%  '
Set-TextValue $ws.Range('B16') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D16') '4.081.83'
Set-TextValue $ws.Range('E16') '  +2.25%  '
Set-TextValue $ws.Range('E17') '  -2.96%  '
Set-TextValue $ws.Range('D18') '8.28'
Set-TextValue $ws.Range('E18') '  -5.89%  '
Set-TextValue $ws.Range('D19') '3.442.12'
Set-TextValue $ws.Range('E19') '  +3.26%  '
Set-TextValue $ws.Range('D20') '17.51'
Set-TextValue $ws.Range('E20') '  -2.82%  '
Set-TextValue $ws.Range('D21') '11.25'
Set-TextValue $ws.Range('E21') '  +3.28%  '
Set-TextValue $ws.Range('D22') '0.485'
Set-TextValue $ws.Range('E22') '  -8.45%  '
Set-TextValue $ws.Range('D23') '494.22'
Set-TextValue $ws.Range('E23') '  -3.77%  '
Set-TextValue $ws.Range('E24') '  -5.29%  '
Set-TextValue $ws.Range('D25') '0.0000192'
Set-TextValue $ws.Range('E25') '  -5.12%  '
Set-TextValue $ws.Range('D26') '6.50'
Set-TextValue $ws.Range('E26') '  -6.59%  '
Set-TextValue $ws.Range('D27') '90.59'
Set-TextValue $ws.Range('E27') '  -6.50%  '
Set-TextValue $ws.Range('D28') '3.621.09'
Set-TextValue $ws.Range('E28') '  +2.14%  '
Set-TextValue $ws.Range('D29') '11.88'
Set-TextValue $ws.Range('E29') '  -4.53%  '
Set-TextValue $ws.Range('D30') '11.62'
Set-TextValue $ws.Range('E30') '  -0.16%  '
Set-TextValue $ws.Range('E31') '  -0.69%  '
Set-TextValue $ws.Range('D32') '2.72'
Set-TextValue $ws.Range('E32') '  +4.15%  '
Set-TextValue $ws.Range('D33') '0.135'
Set-TextValue $ws.Range('E33') '  -7.85%  '
Set-TextValue $ws.Range('E34') '  -6.43%  '
Set-TextValue $ws.Range('E35') '  +0.25%  '
Set-TextValue $ws.Range('D36') '29.91'
Set-TextValue $ws.Range('E36') '  +3.86%  '
Set-TextValue $ws.Range('D37') '0.551'
Set-TextValue $ws.Range('E37') '  -3.44%  '
Set-TextValue $ws.Range('D38') '547.48'
Set-TextValue $ws.Range('E38') '  +4.40%  '
Set-TextValue $ws.Range('B39') 'Fetch.AI'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D39') '1.44'
Set-TextValue $ws.Range('E39') '  -5.70%  '
Set-TextValue $ws.Range('B40') 'RenderToken'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range('D40') '7.55'
Set-TextValue $ws.Range('E40') '  -5.12%  '
Set-TextValue $ws.Range('E41') '  -0.09%  '
Set-TextValue $ws.Range('E42') '  -1.33%  '
Set-TextValue $ws.Range('D43') '0.917'
Set-TextValue $ws.Range('E43') '  +6.41%  '
Set-TextValue $ws.Range('E44') '  -1.71%  '
Set-TextValue $ws.Range('D45') '1.70'
Set-TextValue $ws.Range('E45') '  -2.80%  '
Set-TextValue $ws.Range('D46') '0.0408'
Set-TextValue $ws.Range('E46') '  -7.46%  '
Set-TextValue $ws.Range('D47') '5.49'
Set-TextValue $ws.Range('E47') '  -3.57%  '
Set-TextValue $ws.Range('E48') '  -2.91%  '
Set-TextValue $ws.Range('D49') '2.13'
Set-TextValue $ws.Range('E49') '  +3.48%  '
Set-TextValue $ws.Range('E50') '  +0.62%  '
Set-TextValue $ws.Range('D51') '52.91'
Set-TextValue $ws.Range('E51') '  -1.63%  '
